$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.398.68"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "1.841.71"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.89"
$ws.Range("E5").Value = "  -7.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5105"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3235"
$ws.Range("E8").Value = "  -7.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06724"
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.14"
$ws.Range("E10").Value = "  -4.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7700"
$ws.Range("E11").Value = "  -4.54%  "
$ws.Range("D12").Value = "1.899.93"
$ws.Range("E12").Value = "  +3.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07689"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.96"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.024"
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007878"
$ws.Range("E19").Value = "  -2.29%  "
$ws.Range("D20").Value = "26.460.89"
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("D21").Value = "2.124.54"
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.569"
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.527"
$ws.Range("E23").Value = "  -5.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.976"
$ws.Range("E24").Value = "  -3.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.338"
$ws.Range("E25").Value = "  -1.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.29"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.650"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.92"
$ws.Range("E28").Value = "  -1.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.85"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.209"
$ws.Range("E30").Value = "  -3.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.161"
$ws.Range("E31").Value = "  -3.27%  "
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04810"
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.132"
$ws.Range("E34").Value = "  -3.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.861"
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6807"
$ws.Range("E36").Value = "  -7.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.076"
$ws.Range("E37").Value = "  -4.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01810"
$ws.Range("E38").Value = "  -2.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.225"
$ws.Range("E39").Value = "  -6.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4924"
$ws.Range("E40").Value = "  -4.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "112.93"
$ws.Range("E41").Value = "  -2.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9015"
$ws.Range("E42").Value = "  -6.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.101"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.757"
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4250"
$ws.Range("E46").Value = "  -6.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1284"
$ws.Range("E47").Value = "  -5.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.143"
$ws.Range("E48").Value = "  -2.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05896"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("E50").Value = "  -3.25%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.431"
$ws.Range("E51").Value = "  -4.59%  "
